$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.777099999999997
$ws.Range("C7").Value = -12.7401
$ws.Range("B9").Value = 5.958099999999994
$ws.Range("C12").Value = -10.5877
$ws.Range("B13").Value = 5.878199999999998
$ws.Range("C14").Value = -12.7302
$ws.Range("D15").Value = -8.727599999999999
$ws.Range("B16").Value = 6.246200000000003
$ws.Range("B18").Value = 6.733299999999998
$ws.Range("C19").Value = -12.01500000000001
$ws.Range("B20").Value = 8.9757
$ws.Range("B26").Value = 5.36640000000001
$ws.Range("C26").Value = -12.7176
$ws.Range("B27").Value = 5.603300000000004
$ws.Range("C27").Value = -12.97479999999999
$ws.Range("D28").Value = -8.292900000000001
$ws.Range("B29").Value = 5.189199999999997
$ws.Range("C29").Value = -10.6231
$ws.Range("D33").Value = -7.571800000000001
$ws.Range("B35").Value = 8.259600000000006
$ws.Range("D35").Value = -8.443699999999996
$ws.Range("B36").Value = 9.082000000000006
$ws.Range("C37").Value = -13.4028
$ws.Range("C38").Value = -12.9066
$ws.Range("D38").Value = -8.0604
$ws.Range("D43").Value = -8.2783
$ws.Range("D44").Value = -7.529299999999999
$ws.Range("B45").Value = 5.870500000000001
$ws.Range("D45").Value = -7.819899999999998
$ws.Range("C47").Value = -11.70150000000001
$ws.Range("D47").Value = -7.385300000000001
$ws.Range("C51").Value = -12.3206
$ws.Range("D51").Value = -7.818099999999996
$ws.Range("C52").Value = -11.29800000000001
$ws.Range("D54").Value = -8.135600000000004
$ws.Range("B55").Value = 6.405899999999995
$ws.Range("C55").Value = -13.5407
$ws.Range("B57").Value = 5.160399999999994
$ws.Range("D57").Value = -8.046500000000002
$ws.Range("D62").Value = -8.454799999999997
$ws.Range("D63").Value = -7.976299999999997
$ws.Range("D67").Value = -6.123200000000002
$ws.Range("B69").Value = 6.166099999999992
$ws.Range("C69").Value = -11.2792
$ws.Range("C70").Value = -12.46110000000001
$ws.Range("D70").Value = -8.0665
$ws.Range("B76").Value = 5.474100000000001
$ws.Range("C76").Value = -12.4554
$ws.Range("B78").Value = 9.7805
$ws.Range("C81").Value = -12.46970000000001
$ws.Range("D81").Value = -8.1685
$ws.Range("B82").Value = 5.457300000000004
$ws.Range("B83").Value = 6.328399999999995
$ws.Range("C83").Value = -13.9741
$ws.Range("D88").Value = -7.522999999999995
$ws.Range("B93").Value = 6.137599999999999
$ws.Range("C94").Value = -10.1527
$ws.Range("D96").Value = -7.860100000000004
$ws.Range("B97").Value = 6.160999999999998
$ws.Range("D99").Value = -7.5701
$ws.Range("C100").Value = -12.4945
$ws.Range("C102").Value = -13.2342
